# Update generated output numbers (column F) on the "展览" and "全部类型"
# worksheets, matching the data refresh performed by the gh-pages build.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 130
$ws1.Range("F5").Value  = 6673
$ws1.Range("F9").Value  = 6144
$ws1.Range("F10").Value = 44
$ws1.Range("F12").Value = 1248
$ws1.Range("F13").Value = 1248
$ws1.Range("F15").Value = 93
$ws1.Range("F17").Value = 117
$ws1.Range("F18").Value = 18
$ws1.Range("F19").Value = 360
$ws1.Range("F22").Value = 4471
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 189
$ws1.Range("F26").Value = 48

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 130
$ws4.Range("F5").Value  = 6673
$ws4.Range("F9").Value  = 6145
$ws4.Range("F10").Value = 44
$ws4.Range("F12").Value = 1248
$ws4.Range("F13").Value = 1248
$ws4.Range("F15").Value = 93
$ws4.Range("F17").Value = 117
$ws4.Range("F18").Value = 18
$ws4.Range("F19").Value = 360
$ws4.Range("F22").Value = 4471
$ws4.Range("F25").Value = 27
$ws4.Range("F26").Value = 189
$ws4.Range("F27").Value = 48
